# The author reworked the TestNG listener (onTestStart / onTestSkipped) so that
# the "BankManagerLoginTest" row in the test_suite runmode sheet is now turned
# off ("N" instead of "Y") and the workbook is left with the test_suite sheet
# as the active/selected tab (instead of AddCustomerTest).

$wb = $excel.ActiveWorkbook

$testSuite = $wb.Worksheets.Item("test_suite")

# Flip the Runmode flag for BankManagerLoginTest (row 4, column B) from Y to N
$testSuite.Range("B4").Value = "N"

# Make test_suite the active sheet and move its selection to the edited cell;
# this also clears the "tabSelected" state previously held by AddCustomerTest.
$testSuite.Range("B4").Select()
